$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Summary")

# Day 2 block (rows 8-11): set C10 "Total Execution" row to text "weekend"
$ws.Range("C10").Value = "weekend"

# Day 3 block (rows 14-17)
$ws.Range("C15").Value = 6884
$ws.Range("C16").Value = 1623
$ws.Range("C17").Value = 1020

# Day 4 block (rows 20-23)
$ws.Range("C21").Value = 6936
$ws.Range("C22").Value = 1675
$ws.Range("C23").Value = 1675

# Update view: scroll so that row 5 is the top-left visible cell, and select C23
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("C23").Select()
